$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 11, shifting existing rows 11-22 down to 12-23.
$ws.Rows(11).Insert() | Out-Null

# The insert carries formatting from the row above into column B; remove it
# so the new row only contains the new E11 entry, matching the target layout.
$ws.Range("B11").Clear() | Out-Null

# Add the new reviewer journal entry.
$e11 = $ws.Range("E11")
$e11.Value = "\href{https://royalsocietypublishing.org/journal/rsbl}{Biology Letters}"

# Keep the same formatting as neighboring E-column cells (left/top aligned, wrap text).
$e11.HorizontalAlignment = $ws.Range("E12").HorizontalAlignment
$e11.VerticalAlignment = $ws.Range("E12").VerticalAlignment
$e11.WrapText = $ws.Range("E12").WrapText

# Match the saved cursor position reflected in the workbook.
$ws.Range("E18").Select() | Out-Null
